$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-19 19:12:13", 0.0004),
    @("2023-12-19 19:13:37", 0.003400000000000001),
    @("2023-12-19 19:14:29", 0.0022),
    @("2023-12-19 19:14:35", 0.0004),
    @("2023-12-19 19:14:45", 0.0004)
)

$startRow = 490
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
